$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 25 so the existing rows 25-28 shift down to 26-29.
$ws.Rows.Item(25).Insert()

# Populate the new row 25 with the new weekly record, copying the date cell
# style (s="2") from the row that was pushed down (now row 26).
$ws.Cells.Item(25, 4).NumberFormat = $ws.Cells.Item(26, 4).NumberFormat

$ws.Cells.Item(25, 1).Value = 9
$ws.Cells.Item(25, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(25, 3).Value = "Metropolitana"
$ws.Cells.Item(25, 4).Value = 44782
$ws.Cells.Item(25, 5).Value = 13
$ws.Cells.Item(25, 6).Value = 100112010
$ws.Cells.Item(25, 7).Value = "Achicoria"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 70
$ws.Cells.Item(25, 11).Value = 6000
$ws.Cells.Item(25, 12).Value = 6000
$ws.Cells.Item(25, 13).Value = 6000
$ws.Cells.Item(25, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(25, 15).Value = "Región Metropolitana"
$ws.Cells.Item(25, 16).Value = 375
$ws.Cells.Item(25, 17).Value = 16
$ws.Cells.Item(25, 18).Value = "Hortaliza"
